$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("G1").Value = "Elapsed Time"
$ws.Range("H1").Value = "CPU"

# Copy the header style from an existing header cell (F1) to the new headers
$ws.Range("F1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("G1").Value = "Elapsed Time"
$ws.Range("H1").Value = "CPU"

# New data cells
$ws.Range("G2").Value = 0.1289622459000384
$ws.Range("H2").Value = 0.991
